# "Cleaned up and expanded logging test cases"
#
# The two logging test-case cells (K1, A6) get new sample values. Both
# cells store numeric-looking text ("2854", "8465") as genuine strings
# (not numbers), matching the existing pattern in the sheet where K1/A6
# already held digit-only text ("2407"/"1833") from the shared-string
# table. Setting NumberFormat to Text ("@") before assigning the value
# is what keeps Excel from auto-converting the digit string to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").NumberFormat = "@"
$ws.Range("K1").Value = "2854"

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "8465"
